$d = $word.ActiveDocument

# 1) Remove the centered-paragraph formatting ("jc=center") from the
#    "Offene Stellen" section paragraphs - replaced breaks/centering with
#    normal (left) paragraph formatting.
foreach ($p in $d.Paragraphs) {
    if ($p.Format.Alignment -eq 1) {
        $p.Format.Alignment = 0
    }
}

# 2) Merge the split runs (incl. the spell-check proofErr markers around
#    "HeilerziehungspflegerInnen") into a single run.
$oldList = "(ErzieherInnen, ErzieherInnen im Anerkennungsjahr, KinderpflegerInnen, HeilerziehungspflegerInnen, SPS1 und SPS2 PraktikantInnen, FÖJ, FSJ, …)"
$d.Content.Find.Execute($oldList, $true, $false, $false, $false, $false, $true, 1, $false, $oldList, 2) | Out-Null

# 3) Replace the closing "Wir freuen uns riesig auf weiteres Personal! :-)"
#    paragraph's text smiley with the real emoji run (matching the
#    mc:AlternateContent / w16se:symEx markup already used elsewhere in
#    the document for the inserted emoji symbol).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Wir freuen uns riesig auf weiteres Personal!")) {
        $target = $p
    }
}

$emojiParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex"><w:r><w:t xml:space="preserve">Wir freuen uns riesig auf weiteres Personal! </w:t></w:r><w:r><w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"/><mc:Fallback><w:rFonts w:ascii="Segoe UI Emoji" w:eastAsia="Segoe UI Emoji" w:hAnsi="Segoe UI Emoji" w:cs="Segoe UI Emoji"/></mc:Fallback></mc:AlternateContent></w:rPr><mc:AlternateContent><mc:Choice Requires="w16se"><w16se:symEx w16se:font="Segoe UI Emoji" w16se:char="1F60A"/></mc:Choice><mc:Fallback><w:t>😊</w:t></mc:Fallback></mc:AlternateContent></w:r></w:p>'

$target.Range.InsertXML($emojiParagraphXml) | Out-Null
